# "BSD - Protokoll.xlsx" update
#
# The sheet is a team log with three parallel "task log" tables, one per
# team member, sharing the same row grid:
#   C:G  -> Lamprecht Daniel   (table "Tabelle2",   C6:G39)
#   M:Q  -> Ruhdorfer Alexander(table "Tabelle24",  M6:Q40)
#   W:AA -> Sammer Manuel      (table "Tabelle245", W6:AA40)
# Each block has columns Aufgabe (task) | Datum (date) | Von (from) |
# Bis (to) | Status (percent complete).
#
# This edit appends one new task row for each of the three people:
#   - Ruhdorfer Alexander (row 37): "Google Maps routing einbauen"
#   - Sammer Manuel       (row 32): "App Routing"
#   - Lamprecht Daniel    (row 37): "User Management"
# all dated 2017-03-29, 14:10 - 15:50, with differing completion %.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$taskDate = 42823          # 2017-03-29 (serial date)
$timeFrom = 0.59027777777777779   # 14:10
$timeTo   = 0.65972222222222221   # 15:50

# --- Ruhdorfer Alexander block (M:Q), row 37 ---------------------------
$ws.Range("M37").Value = "Google Maps routing einbauen"
$ws.Range("N37").Value = $taskDate
$ws.Range("O37").Value = $timeFrom
$ws.Range("P37").Value = $timeTo
$ws.Range("Q37").Value = 0.8

# --- Sammer Manuel block (W:AA), row 32 ---------------------------------
$ws.Range("W32").Value = "App Routing"
$ws.Range("X32").Value = $taskDate
$ws.Range("Y32").Value = $timeFrom
$ws.Range("Z32").Value = $timeTo
$ws.Range("AA32").Value = 0.1

# --- Lamprecht Daniel block (C:G), row 37 -------------------------------
$ws.Range("C37").Value = "User Management"
$ws.Range("D37").Value = $taskDate
$ws.Range("E37").Value = $timeFrom
$ws.Range("F37").Value = $timeTo
$ws.Range("G37").Value = 0.7

# Leave the selection where the author's saved view had it.
[void]$ws.Range("G38").Select()
